$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '55.534.49'
$ws.Range("E2").Value = '  -5.80%  '

# Row 3
$ws.Range("D3").Value = '2.931.18'
$ws.Range("E3").Value = '  -9.30%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +1.29%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '475.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -11.57%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '123.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -9.43%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.02'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.56%  '

# Row 8
$ws.Range("D8").Value = '2.926.15'
$ws.Range("E8").Value = '  -9.47%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.398'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -13.29%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.58'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -13.50%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0950'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -17.93%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.330'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -16.82%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.124'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.76%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.476.40'
$ws.Range("E14").Value = '  -8.16%  '

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.026.07'
$ws.Range("E15").Value = '  -6.29%  '

# Row 16
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '56.348.99'
$ws.Range("E16").Value = '  -4.61%  '

# Row 17
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '22.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -14.37%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000131'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -17.79%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -13.90%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -13.24%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -15.45%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '310.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -13.96%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.10%  '

# Row 24
$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.446'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -14.12%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '59.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -16.32%  '

# Row 26
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.16%  '

# Row 27
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.157'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.26%  '

# Row 28
$ws.Range("B28").Value = 'USDe'
$ws.Range("C28").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.07%  '

# Row 29
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0796'
$ws.Range("E29").Value = '  -19.07%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -17.12%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -10.13%  '

# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.11'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -13.24%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -14.85%  '

# Row 34
$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.60'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -17.33%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '145.24'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -10.88%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.13'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -16.44%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -16.47%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.20'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -16.58%  '

# Row 39
$ws.Range("D39").Value = '3.080.40'
$ws.Range("E39").Value = '  -5.59%  '

# Row 40
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.00%  '

# Row 41
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0607'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -14.29%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.95'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -19.61%  '

# Row 43
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '34.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -16.96%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.605'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -16.04%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.930'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -14.98%  '

# Row 46
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.070.60'
$ws.Range("E46").Value = '  -9.79%  '

# Row 47
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -15.90%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -14.71%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.26'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -16.06%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.54'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -15.64%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0210'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -13.54%  '
